$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '72.794.62'
$ws.Range('E2').Value = '  +0.58%  '

# Row 3
$ws.Range('D3').Value = '2.674.81'
$ws.Range('E3').Value = '  +1.97%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.90'
$ws.Range('E5').Value = '  -0.71%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.41'
$ws.Range('E6').Value = '  -1.19%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('E8').Value = '  +0.12%  '

# Row 9
$ws.Range('D9').Value = '2.672.54'
$ws.Range('E9').Value = '  +1.94%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.170'
$ws.Range('E10').Value = '  -1.22%  '

# Row 11
$ws.Range('E11').Value = '  +2.40%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.358'
$ws.Range('E12').Value = '  +1.94%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.02'
$ws.Range('E13').Value = '  -0.33%  '

# Row 14
$ws.Range('D14').Value = '3.160.68'
$ws.Range('E14').Value = '  +1.73%  '

# Row 15
$ws.Range('D15').Value = '72.548.26'
$ws.Range('E15').Value = '  +0.54%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000186'
$ws.Range('E16').Value = '  -0.50%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.38'
$ws.Range('E17').Value = '  -0.61%  '

# Row 18
$ws.Range('D18').Value = '2.664.07'
$ws.Range('E18').Value = '  +1.75%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.27'
$ws.Range('E19').Value = '  +6.14%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.28'
$ws.Range('E20').Value = '  +4.43%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '372.64'
$ws.Range('E21').Value = '  -2.35%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.19'
$ws.Range('E22').Value = '  +0.82%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.05'
$ws.Range('E23').Value = '  +1.42%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.24'
$ws.Range('E24').Value = '  -1.16%  '

# Row 25
$ws.Range('E25').Value = '  +0.05%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.35'
$ws.Range('E26').Value = '  -0.40%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.80'
$ws.Range('E27').Value = '  -1.16%  '

# Row 28
$ws.Range('D28').Value = '2.809.26'
$ws.Range('E28').Value = '  +1.96%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.22%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0977'
$ws.Range('E30').Value = '  +2.80%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.13'
$ws.Range('E31').Value = '  +1.26%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '503.21'
$ws.Range('E32').Value = '  -3.07%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.31'
$ws.Range('E33').Value = '  -1.50%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.83'
$ws.Range('E34').Value = '  +0.24%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.05%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.97'
$ws.Range('E36').Value = '  -0.43%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.61'
$ws.Range('E37').Value = '  +1.81%  '

# Row 38
$ws.Range('B38').Value = 'WhiteBITCoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.97'
$ws.Range('E38').Value = '  -0.59%  '

# Row 39
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.111'
$ws.Range('E39').Value = '  -0.28%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.39'
$ws.Range('E40').Value = '  -0.79%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.79'
$ws.Range('E41').Value = '  -2.08%  '

# Row 42
$ws.Range('E42').Value = '  +0.03%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.02'
$ws.Range('E43').Value = '  -0.41%  '

# Row 44
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.57'
$ws.Range('E44').Value = '  -0.31%  '

# Row 45
$ws.Range('B45').Value = 'PolygonEcosystemToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.334'
$ws.Range('E45').Value = '  +0.84%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '156.96'
$ws.Range('E46').Value = '  +4.52%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.53'
$ws.Range('E47').Value = '  +0.10%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.77'
$ws.Range('E48').Value = '  +2.37%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.561'
$ws.Range('E49').Value = '  +3.52%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.74'
$ws.Range('E50').Value = '  +2.64%  '

# Row 51
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.609'
$ws.Range('E51').Value = '  +1.81%  '
